$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Nikola Jokic, DEN, 1 (rank 1 unchanged)
$ws.Range("B2").Value = "Nikola Jokic"
$ws.Range("C2").Value = "DEN"
$ws.Range("D2").Value = 1

# Row 3: rank 1, Jayson Tatum, BOS, 1
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Jayson Tatum"
$ws.Range("C3").Value = "BOS"
$ws.Range("D3").Value = 1

# Row 4: rank unchanged (3), LeBron James, LAL, 0
$ws.Range("B4").Value = "LeBron James"
$ws.Range("C4").Value = "LAL"
$ws.Range("D4").Value = 0

# Row 5: rank 3, Kyle Lowry, PHI, 0
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Kyle Lowry"
$ws.Range("C5").Value = "PHI"
$ws.Range("D5").Value = 0

# Row 6: rank 3, Mike Conley, MIN, 0
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Mike Conley"
$ws.Range("C6").Value = "MIN"
$ws.Range("D6").Value = 0
